# Scheduled runner update: refresh market price / profit columns (H-N) across Sheets
# Values sourced from a data refresh; applies updates + clears per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 739.6667
$ws.Range("I96").Value = 561.75
$ws.Range("J96").Value = 882
$ws.Range("K96").Value = 1685.25
$ws.Range("L96").Value = 2646
$ws.Range("M96").Value = -312.25
$ws.Range("N96").Value = -5392
$ws.Range("H106").Value = 4782.1875
$ws.Range("I106").Value = 5451.5
$ws.Range("K106").Value = 5451.5
$ws.Range("M106").Value = -4820.5
$ws.Range("H118").Value = 1101.8182
$ws.Range("I118").Value = 270
$ws.Range("J118").Value = 2100
$ws.Range("K118").Value = 810
$ws.Range("L118").Value = 6300
$ws.Range("M118").Value = 847
$ws.Range("N118").Value = -9614
$ws.Range("H123").Value = 85097.5
$ws.Range("J123").Value = 85097.5
$ws.Range("L123").Value = 85097.5
$ws.Range("N123").Value = -94897.5
$ws.Range("H128").Value = 40000
$ws.Range("I128").Value = 40000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 40000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -35020
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 808.5454999999999
$ws.Range("I129").Value = 546.1667
$ws.Range("J129").Value = 1123.4
$ws.Range("K129").Value = 1638.5001
$ws.Range("L129").Value = 3370.2
$ws.Range("M129").Value = 3361.4999
$ws.Range("N129").Value = -13370.2
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H135").Value = 29902.584
$ws.Range("I135").Value = 34015.195
$ws.Range("J135").Value = 4404.4
$ws.Range("K135").Value = 306136.755
$ws.Range("L135").Value = 39639.6
$ws.Range("M135").Value = -303601.755
$ws.Range("N135").Value = -44709.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2271.913
$ws.Range("I2").Value = 1388.7693
$ws.Range("J2").Value = 3420
$ws.Range("K2").Value = 1388.7693
$ws.Range("L2").Value = 3420
$ws.Range("M2").Value = -1275.7693
$ws.Range("N2").Value = -3646
$ws.Range("H116").Value = 2271.913
$ws.Range("I116").Value = 1388.7693
$ws.Range("J116").Value = 3420
$ws.Range("K116").Value = 1388.7693
$ws.Range("L116").Value = 3420
$ws.Range("M116").Value = 905.2307000000001
$ws.Range("N116").Value = -8008
$ws.Range("H122").Value = 9261178
$ws.Range("I122").Value = 2220.4
$ws.Range("J122").Value = 55555964
$ws.Range("K122").Value = 6661.200000000001
$ws.Range("L122").Value = 166667892
$ws.Range("M122").Value = -4211.200000000001
$ws.Range("N122").Value = -166672792
$ws.Range("H132").Value = 162631.47
$ws.Range("I132").Value = 144037.28
$ws.Range("J132").Value = 184324.67
$ws.Range("K132").Value = 432111.84
$ws.Range("L132").Value = 552974.01
$ws.Range("M132").Value = -429581.84
$ws.Range("N132").Value = -558034.01
$ws.Range("H135").Value = 49714.5
$ws.Range("J135").Value = 49714.5
$ws.Range("L135").Value = 49714.5
$ws.Range("N135").Value = -59854.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2271.913
$ws.Range("I3").Value = 1388.7693
$ws.Range("J3").Value = 3420
$ws.Range("K3").Value = 1388.7693
$ws.Range("L3").Value = 3420
$ws.Range("M3").Value = -1274.7693
$ws.Range("N3").Value = -3648
$ws.Range("H94").Value = 1096.3334
$ws.Range("I94").Value = 511.4
$ws.Range("J94").Value = 1827.5
$ws.Range("K94").Value = 511.4
$ws.Range("L94").Value = 1827.5
$ws.Range("M94").Value = -60.39999999999998
$ws.Range("N94").Value = -2729.5
$ws.Range("H99").Value = 1162.3684
$ws.Range("I99").Value = 1055
$ws.Range("J99").Value = 1281.6666
$ws.Range("K99").Value = 1055
$ws.Range("L99").Value = 1281.6666
$ws.Range("M99").Value = 443
$ws.Range("N99").Value = -4277.6666
$ws.Range("H102").Value = 7666.75
$ws.Range("I102").Value = 7666.75
$ws.Range("K102").Value = 7666.75
$ws.Range("M102").Value = -4421.75
$ws.Range("H134").Value = 3993.4
$ws.Range("I134").Value = 3191.75
$ws.Range("J134").Value = 7200
$ws.Range("K134").Value = 9575.25
$ws.Range("L134").Value = 21600
$ws.Range("M134").Value = -7040.25
$ws.Range("N134").Value = -26670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3670.4
$ws.Range("I62").Value = 3266.6667
$ws.Range("J62").Value = 3843.4285
$ws.Range("K62").Value = 3266.6667
$ws.Range("L62").Value = 3843.4285
$ws.Range("M62").Value = -2642.6667
$ws.Range("N62").Value = -5091.4285
$ws.Range("H65").Value = 3670.4
$ws.Range("I65").Value = 3266.6667
$ws.Range("J65").Value = 3843.4285
$ws.Range("K65").Value = 16333.3335
$ws.Range("L65").Value = 19217.1425
$ws.Range("M65").Value = -13213.3335
$ws.Range("N65").Value = -25457.1425

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 457.14285
$ws.Range("H122").Value = 737.7941
$ws.Range("I122").Value = 341.64285
$ws.Range("J122").Value = 1015.1
$ws.Range("K122").Value = 3074.78565
$ws.Range("L122").Value = 9135.9
$ws.Range("M122").Value = -624.7856500000003
$ws.Range("N122").Value = -14035.9
$ws.Range("H131").Value = 862
$ws.Range("I131").Value = 498.7143
$ws.Range("J131").Value = 989.15
$ws.Range("K131").Value = 1496.1429
$ws.Range("L131").Value = 2967.45
$ws.Range("M131").Value = 3543.8571
$ws.Range("N131").Value = -13047.45
$ws.Range("H135").Value = 457.14285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1453.1875
$ws.Range("J122").Value = 2636
$ws.Range("L122").Value = 7908
$ws.Range("N122").Value = -12808
$ws.Range("H126").Value = 2018.5883
$ws.Range("J126").Value = 1956.5
$ws.Range("L126").Value = 5869.5
$ws.Range("N126").Value = -10809.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 988.5217
$ws.Range("I93").Value = 972.625
$ws.Range("J93").Value = 1024.8572
$ws.Range("K93").Value = 972.625
$ws.Range("L93").Value = 1024.8572
$ws.Range("M93").Value = 275.375
$ws.Range("N93").Value = -3520.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1528.091
$ws.Range("I126").Value = 1381
$ws.Range("K126").Value = 4143
$ws.Range("M126").Value = -1673

Write-Host "Applied Hades_Profits scheduled update."
